$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: C6:F6 change from 2 to 5, and G6 gets a new value of 5 (was empty)
$ws.Range("C6:G6").Value = 5

# Update active cell selection to H6
$ws.Range("H6").Select()
